# Commit: "hours set to< 40; charts are good to go."
#
# Changes applied (per the OOXML diff):
#  1. Slide 3 speaker notes: keep the existing "Candida solved our
#     looping list" note, then add a blank paragraph and a new
#     paragraph describing the 800K+ lines / 100 pages / 2000 games
#     detail.
#  2. Slide 4 speaker notes: the old placeholder/draft notes
#     ("We are not adding images...", "Rachel = 2500", "Candida = ?",
#     "Chris = ?", "Max = ?") are removed now that the charts are
#     ready, leaving a single empty paragraph.

$p = $ppt.ActivePresentation

$nl = [char]10

# --- Slide 3 notes: append the new detail paragraphs -----------------
$slide3 = $p.Slides.Item(3)
$notes3 = $slide3.NotesPage.Shapes.Item("Notes Placeholder 2")
$notes3.TextFrame.TextRange.Text = (
    "I really want to highlight Candida solved our looping list for multiple pages of data; I am not sure other would have to have done that. " +
    $nl + $nl +
    "We had over 800K+ lines of data to pull through; we selected 100 pages resulting in 2000 games"
)

# --- Slide 4 notes: clear the old draft notes, charts are good to go --
$slide4 = $p.Slides.Item(4)
$notes4 = $slide4.NotesPage.Shapes.Item("Notes Placeholder 2")
$notes4.TextFrame.TextRange.Text = ""
